$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Price" column (D) values. Set NumberFormat to Text first so that
# numeric-looking strings (e.g. "1.005", "0.000009488") are preserved exactly
# as text instead of being auto-converted to numbers by Excel, then restore
# the original (default) cell style afterwards.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.996.46"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.925.56"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4590"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3820"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07754"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9790"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "22.56"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.932.53"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.710"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.971"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06988"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "84.81"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.006"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009488"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.72"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "29.020.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.346"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.146.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.056"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.18"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.04"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.627"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "117.84"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.842"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09318"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8651"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.108"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.248"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.017"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05699"
$ws.Range("D36").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02054"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.093"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.468"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5512"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1758"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "9.353"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.000002803"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5175"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06940"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "111.06"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.768"
$ws.Range("D51").Style = "Normal"

# Update "Volume(1h)" column (E) values.
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("E3").Value = "  +1.23%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E8").Value = "  +0.30%  "
$ws.Range("E9").Value = "  +0.41%  "
$ws.Range("E10").Value = "  -0.25%  "
$ws.Range("E11").Value = "  +2.42%  "
$ws.Range("E12").Value = "  +1.98%  "
$ws.Range("E13").Value = "  +0.65%  "
$ws.Range("E14").Value = "  +0.09%  "
$ws.Range("E15").Value = "  -1.22%  "
$ws.Range("E16").Value = "  +0.87%  "
$ws.Range("E17").Value = "  -0.05%  "
$ws.Range("E18").Value = "  -0.45%  "
$ws.Range("E19").Value = "  -0.13%  "
$ws.Range("E21").Value = "  +0.61%  "
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("E24").Value = "  +0.56%  "
$ws.Range("E25").Value = "  -2.08%  "
$ws.Range("E26").Value = "  +0.65%  "
$ws.Range("E27").Value = "  -0.64%  "
$ws.Range("E28").Value = "  +0.67%  "
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("E30").Value = "  +0.24%  "
$ws.Range("E31").Value = "  +0.49%  "
$ws.Range("E32").Value = "  +0.38%  "
$ws.Range("E33").Value = "  +0.22%  "
$ws.Range("E34").Value = "  -0.45%  "
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("E36").Value = "  -0.07%  "
$ws.Range("E37").Value = "  +0.50%  "
$ws.Range("E38").Value = "  -0.06%  "
$ws.Range("E39").Value = "  +0.99%  "
$ws.Range("E40").Value = "  +13.45%  "
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("E43").Value = "  +0.32%  "
$ws.Range("E44").Value = "  +0.50%  "
$ws.Range("E45").Value = "  +14.55%  "
$ws.Range("E46").Value = "  +4.68%  "
$ws.Range("E47").Value = "  -0.32%  "
$ws.Range("E48").Value = "  +1.56%  "
$ws.Range("E49").Value = "  -1.09%  "
$ws.Range("E50").Value = "  -0.23%  "
$ws.Range("E51").Value = "  -0.31%  "
